# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values on the zh-cn and
# de-de worksheets to reflect the regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-17 03:30:17"
$wsZhCn.Range("E3").Value = "2016-03-17 03:30:17"
$wsZhCn.Range("H2").Value = "2016-03-17 03:31:03"
$wsZhCn.Range("H3").Value = "2016-03-17 03:31:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-17 03:30:26"
$wsDeDe.Range("E3").Value = "2016-03-17 03:30:26"
$wsDeDe.Range("H2").Value = "2016-03-17 03:31:21"
$wsDeDe.Range("H3").Value = "2016-03-17 03:31:21"
